# Append a new log row (row 9) to the worksheet, recording the timestamp of
# this run together with the same "Random" method statistics columns used
# by the previous rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

# Column A: timestamp (Excel serial date/time) this run happened. Fill down
# from the cell above so it inherits the exact same style (date/time format)
# as every other row in the log, rather than creating a brand new style.
$ws.Range("A8:A9").FillDown()
$ws.Cells.Item($row, 1).Value = 42611.883761574078

# Column B: elapsed milliseconds for this run.
$ws.Cells.Item($row, 2).Value = 76

# Columns C-M: the various score/percentage/count metrics — 0 for this run.
for ($col = 3; $col -le 13; $col++) {
    $ws.Cells.Item($row, $col).Value = 0
}

# Column N: the prediction method used.
$ws.Cells.Item($row, 14).Value = "Random"
